# Apply cryptos list update (prices & 1h volume %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.366.55'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '2.065.62'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.64'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("E6").Value = '  +1.69%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.39'
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.396'
$ws.Range("E9").Value = '  +3.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0775'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.102'
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").Value = '2.369.72'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.40'
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.79'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.775'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.20'
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '2.066.63'
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").Value = '37.324.85'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.29'
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.64'
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").Value = '0.0₃0819'
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.67'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("E25").Value = '  -1.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.19'
$ws.Range("E26").Value = '  +2.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.87'
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.43'
$ws.Range("E28").Value = '  -4.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.12'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.118'
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.53'
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0617'
$ws.Range("E33").Value = '  -0.69%  '
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -3.20%  '
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.64'
$ws.Range("E39").Value = '  -4.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0964'
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.95'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").Value = '1.484.17'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.16'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.07'
$ws.Range("E46").Value = '  -11.71%  '
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.30'
$ws.Range("E48").Value = '  -3.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.23'
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("E50").Value = '  +0.62%  '
$ws.Range("D51").Value = '2.256.66'
$ws.Range("E51").Value = '  +0.13%  '
